$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete rows 23 and 24 (Costsheet and Timesheet) entirely, shifting the
# rows below (including the ArrestWarrant row and the trailing blank rows)
# up by two.
$ws.Rows("23:24").Delete()
